# Apply "MAJ mapping suite review de NRISS" edits to FRCarePlanLMCDAFHIR.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Row 4 = "Name" property: clear its value
$meta.Range("B4").Value = ""

# Row 5 = "Title" property: now carries the value that used to be on "Name"
$meta.Range("B5").Value = 'Mapping Métier/CDA/FHIR : "Référence Item Plan Traitement"'

# Row 8 = "Date" property: bump the timestamp
$meta.Range("B8").Value = "2026-01-07T15:20:53+00:00"

# --- Sheet "Mapping Table 0" ---
$map0 = $wb.Worksheets.Item("Mapping Table 0")

# D6: rename frProduitSante -> consumable
$map0.Range("D6").Value = "FRCDAReferenceItemPlanTraitement.consumable"

# D8: switch entryRelationship "." separator to ":" separator
$map0.Range("D8").Value = "FRCDAReferenceItemPlanTraitement.entryRelationship:frItemPlanTraitement"

# --- Sheet "Mapping Table 1" ---
$map1 = $wb.Worksheets.Item("Mapping Table 1")

# A6: rename frProduitSante -> consumable
$map1.Range("A6").Value = "FRCDAReferenceItemPlanTraitement.consumable"

# A8: switch entryRelationship "." separator to ":" separator
$map1.Range("A8").Value = "FRCDAReferenceItemPlanTraitement.entryRelationship:frItemPlanTraitement"

# D8 & D9: switch reference[...] bracket notation to reference:... colon notation
$map1.Range("D8").Value = "FRCarePlanDocument.activity.reference:FRMedicationRequestDocument"
$map1.Range("D9").Value = "FRCarePlanDocument.activity.reference:FRMedicationRequestDocument"

# D10: same rename, with the .identifier suffix kept
$map1.Range("D10").Value = "FRCarePlanDocument.activity.reference:FRMedicationRequestDocument.identifier"
